$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value that used to be in C109 - it moves down to the new row 110
$ws.Range("C109").Value = ""

# Add the new row 110 with the data previously destined for C109's "NA"
# Force the Date column to stay plain text (matches the source data, which
# stores dates as plain strings rather than Excel date serials), then reset
# the style back to Normal so no extra formatting is left on the cell.
$ws.Range("A110").NumberFormat = "@"
$ws.Range("A110").Value = "2025-05-14"
$ws.Range("A110").Style = "Normal"
$ws.Range("B110").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C110").Value = "NA"
$ws.Range("D110").Value = 1
